$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = -5.7686
$ws.Range("D2").Value = 0.7096
$ws.Range("D2").Interior.Color = 12905163   # 00CBEAC4
$ws.Range("E2").Value = 2.2721
$ws.Range("E2").Interior.Color = 16121079   # 00F7FCF5

# --- Row 3 ---
$ws.Range("C3").Value = -0.5945
$ws.Range("D3").Value = 0.5243
$ws.Range("D3").Interior.Color = 8440448    # 0080CA80
$ws.Range("E3").Value = 1.2607
$ws.Range("E3").Interior.Color = 7914870    # 0076C578

# --- Row 4 ---
$ws.Range("C4").Value = 0.4392
$ws.Range("D4").Value = 0.4152
$ws.Range("D4").Interior.Color = 6467659    # 004BB062
$ws.Range("E4").Value = 1.0114
$ws.Range("E4").Interior.Color = 6204483    # 0043AC5E

# --- Row 5 ---
$ws.Range("C5").Value = 0.9815
$ws.Range("D5").Value = 0.0943
$ws.Range("D5").Interior.Color = 1786880    # 0000441B
$ws.Range("E5").Value = 0.2302
$ws.Range("E5").Interior.Color = 1786880    # 0000441B

# --- Row 6 ---
$ws.Range("C6").Value = 0.6804
$ws.Range("D6").Value = 0.438
$ws.Range("D6").Interior.Color = 6796630    # 0056B567
$ws.Range("E6").Value = 1.0993
$ws.Range("E6").Interior.Color = 6796629    # 0055B567

# --- Row 7 ---
$ws.Range("C7").Value = 0.427
$ws.Range("D7").Value = 0.5888
$ws.Range("D7").Interior.Color = 10016670   # 009ED798
$ws.Range("E7").Value = 1.4735
$ws.Range("E7").Interior.Color = 9885339    # 009BD696

# --- Row 8 ---
$ws.Range("C8").Value = 0.0192
$ws.Range("D8").Value = 0.7832
$ws.Range("D8").Interior.Color = 14480353   # 00E1F3DC
$ws.Range("E8").Value = 1.9327
$ws.Range("E8").Interior.Color = 14086619   # 00DBF1D6

# --- Row 9 ---
$ws.Range("C9").Value = -0.2416
$ws.Range("D9").Value = 0.8987
$ws.Range("D9").Interior.Color = 16121079   # 00F7FCF5
$ws.Range("E9").Value = 2.23
$ws.Range("E9").Interior.Color = 15924212   # 00F4FBF2
